$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 17.29619540223941
$ws.Cells.Item(2, 3).Value = 10.86081673650695
$ws.Cells.Item(2, 4).Value = 5.839207731680963
$ws.Cells.Item(2, 5).Value = 9.411092624575813
$ws.Cells.Item(2, 6).Value = 34.4767854725484
$ws.Cells.Item(2, 13).Value = 16.72178415589129
$ws.Cells.Item(2, 14).Value = 18.65609438729731

$ws.Cells.Item(3, 2).Value = 16.66770173180015
$ws.Cells.Item(3, 3).Value = 10.25130124603323
$ws.Cells.Item(3, 4).Value = 5.855014677497666
$ws.Cells.Item(3, 5).Value = 9.327985292081859
$ws.Cells.Item(3, 6).Value = 33.95775021068366
$ws.Cells.Item(3, 13).Value = 16.43211586361196
$ws.Cells.Item(3, 14).Value = 18.704559713776

$ws.Cells.Item(4, 2).Value = 16.27552129864774
$ws.Cells.Item(4, 3).Value = 9.861107621943347
$ws.Cells.Item(4, 4).Value = 5.865738270268817
$ws.Cells.Item(4, 5).Value = 9.279504279466662
$ws.Cells.Item(4, 6).Value = 33.6467760026637
$ws.Cells.Item(4, 13).Value = 16.25729199141953
$ws.Cells.Item(4, 14).Value = 18.73625885600953

$ws.Cells.Item(5, 2).Value = 16.11440895062061
$ws.Cells.Item(5, 3).Value = 9.698271865777047
$ws.Cells.Item(5, 4).Value = 5.870362626257357
$ws.Cells.Item(5, 5).Value = 9.26040404204689
$ws.Cells.Item(5, 6).Value = 33.52214516453438
$ws.Cells.Item(5, 13).Value = 16.18691059239601
$ws.Cells.Item(5, 14).Value = 18.74966355470369

$ws.Cells.Item(6, 2).Value = 16.08758689441252
$ws.Cells.Item(6, 3).Value = 9.671007875204026
$ws.Cells.Item(6, 4).Value = 5.871145816018217
$ws.Cells.Item(6, 5).Value = 9.257272521522481
$ws.Cells.Item(6, 6).Value = 33.50158115043911
$ws.Cells.Item(6, 13).Value = 16.17527863821537
$ws.Cells.Item(6, 14).Value = 18.75191877925726

$ws.Cells.Item(7, 2).Value = 16.27335333375497
$ws.Cells.Item(7, 3).Value = 9.858926798303038
$ws.Cells.Item(7, 4).Value = 5.865799607965835
$ws.Cells.Item(7, 5).Value = 9.279244010580452
$ws.Cells.Item(7, 6).Value = 33.64508651109998
$ws.Cells.Item(7, 13).Value = 16.25633918962976
$ws.Cells.Item(7, 14).Value = 18.73643766586023

$ws.Cells.Item(8, 2).Value = 17.08095773709585
$ws.Cells.Item(8, 3).Value = 10.65406455213075
$ws.Cells.Item(8, 4).Value = 5.844445704136819
$ws.Cells.Item(8, 5).Value = 9.381918366212053
$ws.Cells.Item(8, 6).Value = 34.29631427241117
$ws.Cells.Item(8, 13).Value = 16.62133595036869
$ws.Cells.Item(8, 14).Value = 18.67240166067969

$ws.Cells.Item(9, 2).Value = 18.6035808152002
$ws.Cells.Item(9, 3).Value = 12.08040843514818
$ws.Cells.Item(9, 4).Value = 5.810720801783428
$ws.Cells.Item(9, 5).Value = 9.602719317891408
$ws.Cells.Item(9, 6).Value = 35.62764776307572
$ws.Cells.Item(9, 13).Value = 17.35683165452813
$ws.Cells.Item(9, 14).Value = 18.56227980556226

$ws.Cells.Item(10, 2).Value = 19.67171448310215
$ws.Cells.Item(10, 3).Value = 13.04064062650098
$ws.Cells.Item(10, 4).Value = 5.791011652195978
$ws.Cells.Item(10, 5).Value = 9.775709390113901
$ws.Cells.Item(10, 6).Value = 36.62942440421061
$ws.Cells.Item(10, 13).Value = 17.90346753237168
$ws.Cells.Item(10, 14).Value = 18.49086276867227

$ws.Cells.Item(11, 2).Value = 20.14433384215461
$ws.Cells.Item(11, 3).Value = 13.45745444288348
$ws.Cells.Item(11, 4).Value = 5.783167833827723
$ws.Cells.Item(11, 5).Value = 9.856491597936319
$ws.Cells.Item(11, 6).Value = 37.08828772162371
$ws.Cells.Item(11, 13).Value = 18.15231892267443
$ws.Cells.Item(11, 14).Value = 18.46044964086801

$ws.Cells.Item(12, 2).Value = 20.32122904065501
$ws.Cells.Item(12, 3).Value = 13.61235513597237
$ws.Cells.Item(12, 4).Value = 5.780360740087708
$ws.Cells.Item(12, 5).Value = 9.88735957965496
$ws.Cells.Item(12, 6).Value = 37.26233564311804
$ws.Cells.Item(12, 13).Value = 18.24648372154766
$ws.Cells.Item(12, 14).Value = 18.44923278991236

$ws.Cells.Item(13, 2).Value = 20.28322622476487
$ws.Cells.Item(13, 3).Value = 13.57912599240406
$ws.Cells.Item(13, 4).Value = 5.780958010379282
$ws.Cells.Item(13, 5).Value = 9.880699632615105
$ws.Cells.Item(13, 6).Value = 37.22484118416133
$ws.Cells.Item(13, 13).Value = 18.22620824482187
$ws.Cells.Item(13, 14).Value = 18.45163517374446

$ws.Cells.Item(14, 2).Value = 20.15892958924
$ws.Cells.Item(14, 3).Value = 13.4702573962845
$ws.Cells.Item(14, 4).Value = 5.782933611309655
$ws.Cells.Item(14, 5).Value = 9.859025695714351
$ws.Cells.Item(14, 6).Value = 37.10260195040454
$ws.Cells.Item(14, 13).Value = 18.16006785458014
$ws.Cells.Item(14, 14).Value = 18.45952080329769

$ws.Cells.Item(15, 2).Value = 20.08251945131359
$ws.Cells.Item(15, 3).Value = 13.40318807055227
$ws.Cells.Item(15, 4).Value = 5.784165031281026
$ws.Cells.Item(15, 5).Value = 9.845785239874836
$ws.Cells.Item(15, 6).Value = 37.02775930226772
$ws.Cells.Item(15, 13).Value = 18.11954307031191
$ws.Cells.Item(15, 14).Value = 18.46439008800884

$ws.Cells.Item(16, 2).Value = 19.64054581475943
$ws.Cells.Item(16, 3).Value = 13.01299282073057
$ws.Cells.Item(16, 4).Value = 5.791546999089976
$ws.Cells.Item(16, 5).Value = 9.770470168145737
$ws.Cells.Item(16, 6).Value = 36.59948566380625
$ws.Cells.Item(16, 13).Value = 17.88719989589468
$ws.Cells.Item(16, 14).Value = 18.49289222311612

$ws.Cells.Item(17, 2).Value = 19.36588329044434
$ws.Cells.Item(17, 3).Value = 12.76844976273579
$ws.Cells.Item(17, 4).Value = 5.796364274225985
$ws.Cells.Item(17, 5).Value = 9.724785616391513
$ws.Cells.Item(17, 6).Value = 36.33744074762073
$ws.Cells.Item(17, 13).Value = 17.74464234408299
$ws.Cells.Item(17, 14).Value = 18.510909884177

$ws.Cells.Item(18, 2).Value = 19.206663224379
$ws.Cells.Item(18, 3).Value = 12.62591522515745
$ws.Cells.Item(18, 4).Value = 5.799240537232763
$ws.Cells.Item(18, 5).Value = 9.698707019682169
$ws.Cells.Item(18, 6).Value = 36.18702571638374
$ws.Cells.Item(18, 13).Value = 17.6626698723302
$ws.Cells.Item(18, 14).Value = 18.52146836234491

$ws.Cells.Item(19, 2).Value = 19.15254622320106
$ws.Cells.Item(19, 3).Value = 12.57733471984115
$ws.Cells.Item(19, 4).Value = 5.800232452911377
$ws.Cells.Item(19, 5).Value = 9.68991193764349
$ws.Cells.Item(19, 6).Value = 36.13615551252168
$ws.Cells.Item(19, 13).Value = 17.63492232240069
$ws.Cells.Item(19, 14).Value = 18.5250767593743

$ws.Cells.Item(20, 2).Value = 19.3952513375291
$ws.Cells.Item(20, 3).Value = 12.79467683423343
$ws.Cells.Item(20, 4).Value = 5.79584053549775
$ws.Cells.Item(20, 5).Value = 9.72962849117361
$ws.Cells.Item(20, 6).Value = 36.36530539749585
$ws.Cells.Item(20, 13).Value = 17.75981618209864
$ws.Cells.Item(20, 14).Value = 18.50897166007405

$ws.Cells.Item(21, 2).Value = 20.19549607187578
$ws.Cells.Item(21, 3).Value = 13.50231489355206
$ws.Cells.Item(21, 4).Value = 5.78234888566513
$ws.Cells.Item(21, 5).Value = 9.865384508578812
$ws.Cells.Item(21, 6).Value = 37.13850011338286
$ws.Cells.Item(21, 13).Value = 18.17949753376443
$ws.Cells.Item(21, 14).Value = 18.45719644957579

$ws.Cells.Item(22, 2).Value = 20.70634198437668
$ws.Cells.Item(22, 3).Value = 13.94765479211746
$ws.Cells.Item(22, 4).Value = 5.77448337824492
$ws.Cells.Item(22, 5).Value = 9.955715605511841
$ws.Cells.Item(22, 6).Value = 37.64542234166879
$ws.Cells.Item(22, 13).Value = 18.45333472478561
$ws.Cells.Item(22, 14).Value = 18.42510741093436

$ws.Cells.Item(23, 2).Value = 20.43485610602167
$ws.Cells.Item(23, 3).Value = 13.71155403406797
$ws.Cells.Item(23, 4).Value = 5.77859361798466
$ws.Cells.Item(23, 5).Value = 9.907364824328567
$ws.Cells.Item(23, 6).Value = 37.3747764499812
$ws.Cells.Item(23, 13).Value = 18.3072545278696
$ws.Cells.Item(23, 14).Value = 18.44207335981085

$ws.Cells.Item(24, 2).Value = 19.38197812009776
$ws.Cells.Item(24, 3).Value = 12.78282562389478
$ws.Cells.Item(24, 4).Value = 5.796076985422631
$ws.Cells.Item(24, 5).Value = 9.727438446226765
$ws.Cells.Item(24, 6).Value = 36.35270703503122
$ws.Cells.Item(24, 13).Value = 17.7529561269282
$ws.Cells.Item(24, 14).Value = 18.50984730871378

$ws.Cells.Item(25, 2).Value = 18.19967843479023
$ws.Cells.Item(25, 3).Value = 11.70955322661461
$ws.Cells.Item(25, 4).Value = 5.818961487626517
$ws.Cells.Item(25, 5).Value = 9.541012466094122
$ws.Cells.Item(25, 6).Value = 35.26264829729023
$ws.Cells.Item(25, 13).Value = 17.15634805218525
$ws.Cells.Item(25, 14).Value = 18.59040941441016

